$wb = $excel.ActiveWorkbook

# --- KMI30 (sheet2): add KSE100, KSE30, ALLSHR ---
# Write order matters for the shared-string table indices that Excel
# assigns to brand-new unique strings: KSE100 -> 267, KSE30 -> 268,
# ALLSHR -> 269 (QSE's GNRI becomes 270 further down).
$wsKMI30 = $wb.Worksheets.Item("KMI30")
$wsKMI30.Range("A31").Value = "KSE100"
$wsKMI30.Range("A33").Value = "KSE30"
$wsKMI30.Range("A32").Value = "ALLSHR"
$wsKMI30.Activate()
$wsKMI30.Range("A31:A33").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1

# --- KMI100 (sheet3): add KSE100, ALLSHR, KSE30 ---
$wsKMI100 = $wb.Worksheets.Item("KMI100")
$wsKMI100.Range("A65").Value = "KSE100"
$wsKMI100.Range("A66").Value = "ALLSHR"
$wsKMI100.Range("A67").Value = "KSE30"
$wsKMI100.Activate()
$wsKMI100.Range("A65:A67").Select()
$excel.ActiveWindow.ScrollRow = 55
$excel.ActiveWindow.ScrollColumn = 1

# --- KMIALL (sheet4): add KSE100, ALLSHR, KSE30 ---
$wsKMIALL = $wb.Worksheets.Item("KMIALL")
$wsKMIALL.Range("A212").Value = "KSE100"
$wsKMIALL.Range("A213").Value = "ALLSHR"
$wsKMIALL.Range("A214").Value = "KSE30"
$wsKMIALL.Activate()
$wsKMIALL.Range("A212:A214").Select()
$excel.ActiveWindow.ScrollRow = 204
$excel.ActiveWindow.ScrollColumn = 1

# --- QSE (sheet5): add GNRI ---
$wsQSE = $wb.Worksheets.Item("QSE")
$wsQSE.Range("A54").Value = "GNRI"
$wsQSE.Activate()
$wsQSE.Range("A54").Select()
$excel.ActiveWindow.ScrollRow = 49
$excel.ActiveWindow.ScrollColumn = 1

# KMIALL ends up as the active sheet/tab, with its new rows selected.
$wsKMIALL.Activate()
$wsKMIALL.Range("A212:A214").Select()
